$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the runs (C) and balls (D) values among rows 4, 5, 6
$ws.Range("C4").Value = "0"
$ws.Range("D4").Value = "0"

$ws.Range("C5").Value = "5"
$ws.Range("D5").Value = "6"

$ws.Range("C6").Value = "7"
$ws.Range("D6").Value = "12"
